# Commit: "finish validation TODO : import graphic, complete TEXT update timetable"
#
# Marks the "done?" (column D) cells with "x" for the rows that were
# missing them (rows 18-46 on Tabelle1 / sheet1), mirroring the pattern
# already used for rows 2-16 directly above them. Also moves the sheet's
# current selection down to where the user was last working (C48).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

for ($r = 18; $r -le 46; $r++) {
    $ws.Cells.Item($r, 4).Value = "x"
}

# Leave the view/selection where the editing session ended up.
$ws.Range("C48").Select()
